# Update Madigan bike hours: Riders (column C) and Average (column D)
# for the Ridership worksheet, rows 2-8 (15 May 2017 - 21 May 2017).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Row -> (Riders, Average)
$updates = @(
    @{ Row = 2; Riders = 249; Average = 226.96 },
    @{ Row = 3; Riders = 183; Average = 219.66 },
    @{ Row = 4; Riders = 243; Average = 214.69 },
    @{ Row = 5; Riders = 244; Average = 234.89 },
    @{ Row = 6; Riders = 221; Average = 238.97 },
    @{ Row = 7; Riders = 154; Average = 112.87 },
    @{ Row = 8; Riders = 72;  Average = 90.5 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.Riders
    $ws.Cells.Item($u.Row, 4).Value = $u.Average
}

$excel.CalculateFullRebuild()
